# Create new worksheet "usaCities" after the last existing sheet (myTrips)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "usaCities"

# Header
$ws.Range("A1").Value = "US Cities"
$ws.Range("A1").Font.Bold = $true

# US city listing (method to navigate to any home listing in the US)
$cities = @(
    "Austin",
    "Tampa",
    "New York",
    "Tysons",
    "Minneapolis",
    "Cincinnati",
    "San Francisco",
    "Sacramento",
    "Nashville",
    "Atlanta",
    "Boston",
    "Portland",
    "Houston",
    "Dallas",
    "Seattle"
)

$row = 2
foreach ($city in $cities) {
    $ws.Cells.Item($row, 1).Value = $city
    $row++
}
